$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23902.5
$ws.Range("J21").Value = 17667.334
$ws.Range("L21").Value = 17667.334
$ws.Range("N21").Value = -18603.334

$ws.Range("H23").Value = 23902.5
$ws.Range("J23").Value = 17667.334
$ws.Range("L23").Value = 17667.334
$ws.Range("N23").Value = -18135.334

$ws.Range("M29").ClearContents()
$ws.Range("H29").Value = 39
$ws.Range("I29").Value = 39
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 117
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = 164

$ws.Range("M38").ClearContents()
$ws.Range("H38").Value = 60.2
$ws.Range("I38").Value = 60.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 180.6
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = 191.4

$ws.Range("H58").Value = 2272
$ws.Range("I58").Value = 1865
$ws.Range("J58").Value = 3900
$ws.Range("K58").Value = 5595
$ws.Range("L58").Value = 11700
$ws.Range("M58").Value = -5445
$ws.Range("N58").Value = -12000

$ws.Range("H76").Value = 3284.4614
$ws.Range("I76").Value = 3370
$ws.Range("J76").Value = 2999.3333
$ws.Range("K76").Value = 3370
$ws.Range("L76").Value = 2999.3333
$ws.Range("M76").Value = -3055
$ws.Range("N76").Value = -3629.3333

$ws.Range("H79").Value = 3284.4614
$ws.Range("I79").Value = 3370
$ws.Range("J79").Value = 2999.3333
$ws.Range("K79").Value = 3370
$ws.Range("L79").Value = 2999.3333
$ws.Range("M79").Value = -2278
$ws.Range("N79").Value = -5183.3333

$ws.Range("H86").Value = 8001.375
$ws.Range("I86").Value = 7000.75
$ws.Range("J86").Value = 9002
$ws.Range("K86").Value = 7000.75
$ws.Range("L86").Value = 9002
$ws.Range("M86").Value = -5877.75
$ws.Range("N86").Value = -11248

$ws.Range("H87").Value = 27530.8
$ws.Range("I87").Value = 23000
$ws.Range("J87").Value = 28663.5
$ws.Range("K87").Value = 23000
$ws.Range("L87").Value = 28663.5
$ws.Range("M87").Value = -21752
$ws.Range("N87").Value = -31159.5

$ws.Range("H89").Value = 8001.375
$ws.Range("I89").Value = 7000.75
$ws.Range("J89").Value = 9002
$ws.Range("K89").Value = 35003.75
$ws.Range("L89").Value = 45010
$ws.Range("M89").Value = -29387.75
$ws.Range("N89").Value = -56242

$ws.Range("H90").Value = 27530.8
$ws.Range("I90").Value = 23000
$ws.Range("J90").Value = 28663.5
$ws.Range("K90").Value = 69000
$ws.Range("L90").Value = 85990.5
$ws.Range("M90").Value = -62760
$ws.Range("N90").Value = -98470.5

$ws.Range("H132").Value = 5066.069
$ws.Range("I132").Value = 1419.9231
$ws.Range("J132").Value = 36666
$ws.Range("K132").Value = 4259.7693
$ws.Range("L132").Value = 109998
$ws.Range("M132").Value = -1729.7693
$ws.Range("N132").Value = -115058

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7576749.5
$ws.Range("I45").Value = 8265454
$ws.Range("K45").Value = 8265454
$ws.Range("M45").Value = -8265077

$ws.Range("H63").Value = 4087.7778
$ws.Range("I63").Value = 2947.5
$ws.Range("K63").Value = 2947.5
$ws.Range("M63").Value = -2261.5

$ws.Range("H66").Value = 4087.7778
$ws.Range("I66").Value = 2947.5
$ws.Range("K66").Value = 14737.5
$ws.Range("M66").Value = -11305.5

$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -26996

$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -84984

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N82").ClearContents()
$ws.Range("H82").Value = 4942.75
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0

$ws.Range("N85").ClearContents()
$ws.Range("H85").Value = 4942.75
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0

$ws.Range("H86").Value = 2052.7273
$ws.Range("I86").Value = 2211.4285
$ws.Range("K86").Value = 2211.4285
$ws.Range("M86").Value = -1088.4285

$ws.Range("H89").Value = 2052.7273
$ws.Range("I89").Value = 2211.4285
$ws.Range("K89").Value = 11057.1425
$ws.Range("M89").Value = -5441.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2300.682
$ws.Range("I31").Value = 1042.3846
$ws.Range("J31").Value = 4118.222
$ws.Range("K31").Value = 1042.3846
$ws.Range("L31").Value = 4118.222
$ws.Range("M31").Value = -747.3846000000001
$ws.Range("N31").Value = -4708.222

$ws.Range("H34").Value = 2300.682
$ws.Range("I34").Value = 1042.3846
$ws.Range("J34").Value = 4118.222
$ws.Range("K34").Value = 1042.3846
$ws.Range("L34").Value = 4118.222
$ws.Range("M34").Value = -840.3846000000001
$ws.Range("N34").Value = -4522.222

$ws.Range("H141").Value = 47573.168
$ws.Range("J141").Value = 50079.816
$ws.Range("L141").Value = 50079.816
$ws.Range("N141").Value = -60439.816

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 463055.84
$ws.Range("I2").Value = 844.9167
$ws.Range("J2").Value = 1255417.4
$ws.Range("K2").Value = 5069.5002
$ws.Range("L2").Value = 7532504.399999999
$ws.Range("M2").Value = -4956.5002
$ws.Range("N2").Value = -7532730.399999999

$ws.Range("H12").Value = 242.06451
$ws.Range("I12").Value = 266.29413
$ws.Range("J12").Value = 212.64285
$ws.Range("K12").Value = 798.88239
$ws.Range("L12").Value = 637.9285500000001
$ws.Range("M12").Value = -625.88239
$ws.Range("N12").Value = -983.9285500000001

$ws.Range("M86").ClearContents()
$ws.Range("H86").Value = 860
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 860
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2580
$ws.Range("N86").Value = -4952

$ws.Range("M89").ClearContents()
$ws.Range("H89").Value = 860
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 860
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 7740
$ws.Range("N89").Value = -19596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5219.737
$ws.Range("I70").Value = 4774.1816
$ws.Range("J70").Value = 5832.375
$ws.Range("K70").Value = 4774.1816
$ws.Range("L70").Value = 5832.375
$ws.Range("M70").Value = -4504.1816
$ws.Range("N70").Value = -6372.375

$ws.Range("H73").Value = 5219.737
$ws.Range("I73").Value = 4774.1816
$ws.Range("J73").Value = 5832.375
$ws.Range("K73").Value = 4774.1816
$ws.Range("L73").Value = 5832.375
$ws.Range("M73").Value = -3838.1816
$ws.Range("N73").Value = -7704.375

$ws.Range("H102").Value = 6537273.5
$ws.Range("I102").Value = 10101985
$ws.Range("J102").Value = 1969
$ws.Range("K102").Value = 10101985
$ws.Range("L102").Value = 1969
$ws.Range("M102").Value = -10100363
$ws.Range("N102").Value = -5213

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4038.3333
$ws.Range("I40").Value = 4229.143
$ws.Range("J40").Value = 3871.375
$ws.Range("K40").Value = 4229.143
$ws.Range("L40").Value = 3871.375
$ws.Range("M40").Value = -4093.143
$ws.Range("N40").Value = -4143.375

$ws.Range("H57").Value = 18015.334
$ws.Range("J57").Value = 18015.334
$ws.Range("L57").Value = 18015.334
$ws.Range("N57").Value = -19147.334

$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 4300
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4300
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12900
$ws.Range("N122").Value = -17800

$ws.Range("H132").Value = 5238.3184
$ws.Range("I132").Value = 5132.067
$ws.Range("J132").Value = 5466
$ws.Range("K132").Value = 15396.201
$ws.Range("L132").Value = 16398
$ws.Range("M132").Value = -12866.201
$ws.Range("N132").Value = -21458

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 71735
$ws.Range("I122").Value = 684.375
$ws.Range("J122").Value = 134891.11
$ws.Range("K122").Value = 2053.125
$ws.Range("L122").Value = 404673.33
$ws.Range("M122").Value = 396.875
$ws.Range("N122").Value = -409573.33

$ws.Range("H126").Value = 933.4
$ws.Range("I126").Value = 979.25
$ws.Range("J126").Value = 750
$ws.Range("K126").Value = 2937.75
$ws.Range("L126").Value = 2250
$ws.Range("M126").Value = -467.75
